# Shared strings are appended in first-write order, so write new unique
# strings in the same order as the target workbook:
#   18 = "Keycaps Filament"
#   19 = "1 Stück von Felix gespendet"
#   20 = "Gesamtkosten pro Person (ohne Versandkst)"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (new: Keycaps Filament) -- written first so its string gets index 18
$ws.Range("A9").Value = "Keycaps Filament"

# Row 3 (ESP32): Stück, Kosten pro Stück, Kosten Gesamt, new "Sonstiges" note
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 35.49
$ws.Range("F3").Value = "1 Stück von Felix gespendet"

# Row 4 (WD2812): Stück, Kosten Gesamt
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 24.2

# Row 5 (innengewinde): Stück, Kosten Gesamt
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 7.55

# Row 8 (Switches option 2): Kosten pro Stück, Stück, Kosten Gesamt (formula)
$ws.Range("B8").Value = 0.37
$ws.Range("C8").Value = 72
$ws.Range("D8").Formula = "=C8*B8"

# Row 9 (new: Keycaps Filament) remaining cells
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 7.9164000000000003

# Row 14 (new: Gesamtkosten pro Person)
$ws.Range("A14").Value = "Gesamtkosten pro Person (ohne Versandkst)"
$ws.Range("D14").Formula = "=(D9+D8+D5+D4+D3)/5"

# New column F width (bestFit-style autosize to fit the longest F entry;
# the COM layer quantizes ColumnWidth to 1/6-character steps, so 23 is the
# closest settable value to the target stored width of ~23.78)
$ws.Columns.Item(6).ColumnWidth = 23

# Update selection to match target
$ws.Range("D10").Select()
